$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.790.65"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "2.665.27"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.25"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.50"
$ws.Range("E6").Value = "  +0.53%  "
$ws.Range("E7").Value = "  +5.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -2.41%  "
$ws.Range("E10").Value = "  +0.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.84"
$ws.Range("E11").Value = "  +0.10%  "
$ws.Range("E12").Value = "  +1.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.94"
$ws.Range("E13").Value = "  -1.06%  "
$ws.Range("E14").Value = "  -2.69%  "
$ws.Range("D15").Value = "3.143.80"
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("D16").Value = "65.702.21"
$ws.Range("E16").Value = "  +0.50%  "
$ws.Range("D17").Value = "2.676.03"
$ws.Range("E17").Value = "  +0.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.56"
$ws.Range("E18").Value = "  -1.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.79"
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "349.45"
$ws.Range("E20").Value = "  -0.32%  "
$ws.Range("E21").Value = "  -1.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.70"
$ws.Range("E23").Value = "  +0.29%  "
$ws.Range("E24").Value = "  +10.90%  "
$ws.Range("E25").Value = "  +0.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.59"
$ws.Range("E26").Value = "  -0.67%  "
$ws.Range("E27").Value = "  +2.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "568.89"
$ws.Range("E28").Value = "  +7.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.17"
$ws.Range("E29").Value = "  +0.78%  "
$ws.Range("E30").Value = "  -2.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.995"
$ws.Range("E31").Value = "  -0.62%  "
$ws.Range("E32").Value = "  -0.48%  "
$ws.Range("E33").Value = "  +3.55%  "
$ws.Range("E34").Value = "  +4.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.51"
$ws.Range("E35").Value = "  +0.69%  "
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.56"
$ws.Range("E37").Value = "  +0.73%  "
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("E39").Value = "  +0.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "154.55"
$ws.Range("E40").Value = "  -2.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "160.48"
$ws.Range("E41").Value = "  -2.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.09"
$ws.Range("E42").Value = "  -0.91%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0615"
$ws.Range("E43").Value = "  +1.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.29"
$ws.Range("E44").Value = "  -1.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.89"
$ws.Range("E45").Value = "  +0.48%  "
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("E47").Value = "  -1.13%  "
$ws.Range("E48").Value = "  +1.03%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.86"
$ws.Range("E49").Value = "  -0.62%  "
$ws.Range("E50").Value = "  -4.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.808"
$ws.Range("E51").Value = "  -1.12%  "
